$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# 1) Queue lock: "Save" column (E) flips from TRUE to FALSE for rows 44-67
for ($r = 44; $r -le 67; $r++) {
    $ws.Cells.Item($r, 5).Value = $false
}

# 2) Player data now saved without the special highlighted formatting on
#    rows 76-77 (GameID / GateID). Clear the highlighted style, then
#    restore row 76's column A to the plain "alt font" style used
#    elsewhere in the sheet (copied from A68, which already carries it).
$ws.Range("A76:J76").ClearFormats()
$ws.Range("A77:J77").ClearFormats()

$ws.Range("A68").Copy()
$ws.Range("A76").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Update the active selection to H78
$ws.Range("H78").Select()
